$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Inhbb"
$ws.Range("C2").Value = "Acvr2a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.4201923333333333
$ws.Range("H2").Value = 1.260577
$ws.Range("I2").Value = 0.08716480679187069
$ws.Range("J2").Value = 0.08716480679187069
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.61878266666667
$ws.Range("N2").Value = 43.856348
$ws.Range("O2").Value = 0.2662829816142094
$ws.Range("P2").Value = 0.2662829816142094
$ws.Range("Q2").Value = 6.142700399199556
$ws.Range("R2").Value = 55.284303592796
$ws.Range("S2").Value = 0.02321050464436582
$ws.Range("T2").Value = 0.02321050464436582

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Inhbb"
$ws.Range("C3").Value = "Acvr2a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.4201923333333333
$ws.Range("H3").Value = 1.260577
$ws.Range("I3").Value = 0.08716480679187069
$ws.Range("J3").Value = 0.08716480679187069
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 27.084169
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.4933415757187404
$ws.Range("P3").Value = 0.4933415757187404
$ws.Range("Q3").Value = 11.38056016850433
$ws.Range("R3").Value = 102.425041516539
$ws.Range("S3").Value = 0.04300202312992105
$ws.Range("T3").Value = 0.04300202312992105

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Inhbb"
$ws.Range("C4").Value = "Acvr2a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.4201923333333333
$ws.Range("H4").Value = 1.260577
$ws.Range("I4").Value = 0.08716480679187069
$ws.Range("J4").Value = 0.08716480679187069
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.19647366666667
$ws.Range("N4").Value = 39.589421
$ws.Range("O4").Value = 0.2403754426670501
$ws.Range("P4").Value = 0.2403754426670501
$ws.Range("Q4").Value = 5.545057061768556
$ws.Range("R4").Value = 49.905513555917
$ws.Range("S4").Value = 0.02095227901758382
$ws.Range("T4").Value = 0.02095227901758381

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Inhbb"
$ws.Range("C5").Value = "Acvr2a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.446732
$ws.Range("H5").Value = 7.340196000000001
$ws.Range("I5").Value = 0.5075507217365239
$ws.Range("J5").Value = 0.5075507217365239
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 14.61878266666667
$ws.Range("N5").Value = 43.856348
$ws.Range("O5").Value = 0.2662829816142094
$ws.Range("P5").Value = 0.2662829816142094
$ws.Range("Q5").Value = 35.76824335157867
$ws.Range("R5").Value = 321.914190164208
$ws.Range("S5").Value = 0.1351521195044455
$ws.Range("T5").Value = 0.1351521195044455

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Inhbb"
$ws.Range("C6").Value = "Acvr2a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.446732
$ws.Range("H6").Value = 7.340196000000001
$ws.Range("I6").Value = 0.5075507217365239
$ws.Range("J6").Value = 0.5075507217365239
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 27.084169
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.4933415757187404
$ws.Range("P6").Value = 0.4933415757187404
$ws.Range("Q6").Value = 66.26770298570801
$ws.Range("R6").Value = 596.4093268713721
$ws.Range("S6").Value = 0.2503958728186806
$ws.Range("T6").Value = 0.2503958728186806

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Inhbb"
$ws.Range("C7").Value = "Acvr2a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.446732
$ws.Range("H7").Value = 7.340196000000001
$ws.Range("I7").Value = 0.5075507217365239
$ws.Range("J7").Value = 0.5075507217365239
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.19647366666667
$ws.Range("N7").Value = 39.589421
$ws.Range("O7").Value = 0.2403754426670501
$ws.Range("P7").Value = 0.2403754426670501
$ws.Range("Q7").Value = 32.28823440739068
$ws.Range("R7").Value = 290.594109666516
$ws.Range("S7").Value = 0.1220027294133977
$ws.Range("T7").Value = 0.1220027294133977

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Inhbb"
$ws.Range("C8").Value = "Acvr2a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.953740666666667
$ws.Range("H8").Value = 5.861222
$ws.Range("I8").Value = 0.4052844714716054
$ws.Range("J8").Value = 0.4052844714716054
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.61878266666667
$ws.Range("N8").Value = 43.856348
$ws.Range("O8").Value = 0.2662829816142094
$ws.Range("P8").Value = 0.2662829816142094
$ws.Range("Q8").Value = 28.56131019302844
$ws.Range("R8").Value = 257.051791737256
$ws.Range("S8").Value = 0.1079203574653981
$ws.Range("T8").Value = 0.1079203574653981

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Inhbb"
$ws.Range("C9").Value = "Acvr2a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.953740666666667
$ws.Range("H9").Value = 5.861222
$ws.Range("I9").Value = 0.4052844714716054
$ws.Range("J9").Value = 0.4052844714716054
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.084169
$ws.Range("N9").Value = 81.25250700000001
$ws.Range("O9").Value = 0.4933415757187404
$ws.Range("P9").Value = 0.4933415757187404
$ws.Range("Q9").Value = 52.91544239817267
$ws.Range("R9").Value = 476.2389815835541
$ws.Range("S9").Value = 0.1999436797701387
$ws.Range("T9").Value = 0.1999436797701387

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Inhbb"
$ws.Range("C10").Value = "Acvr2a"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.953740666666667
$ws.Range("H10").Value = 5.861222
$ws.Range("I10").Value = 0.4052844714716054
$ws.Range("J10").Value = 0.4052844714716054
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.19647366666667
$ws.Range("N10").Value = 39.589421
$ws.Range("O10").Value = 0.2403754426670501
$ws.Range("P10").Value = 0.2403754426670501
$ws.Range("Q10").Value = 25.78248725916244
$ws.Range("R10").Value = 232.042385332462
$ws.Range("S10").Value = 0.0974204342360686
$ws.Range("T10").Value = 0.09742043423606858
